$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Fgf13"
$ws.Range("C2").Value = "Scn8a"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.302256
$ws.Range("H2").Value = 0.906768
$ws.Range("I2").Value = 0.2472511157095891
$ws.Range("J2").Value = 0.2472511157095891
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.2569
$ws.Range("N2").Value = 0.7706999999999999
$ws.Range("O2").Value = 0.1883407151228772
$ws.Range("P2").Value = 0.1883407151228772
$ws.Range("Q2").Value = 0.07764956639999999
$ws.Range("R2").Value = 0.6988460975999999
$ws.Range("S2").Value = 0.04656745194767328
$ws.Range("T2").Value = 0.04656745194767326

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Fgf13"
$ws.Range("C3").Value = "Scn8a"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.302256
$ws.Range("H3").Value = 0.906768
$ws.Range("I3").Value = 0.2472511157095891
$ws.Range("J3").Value = 0.2472511157095891
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.965375
$ws.Range("N3").Value = 2.896125
$ws.Range("O3").Value = 0.7077439387378264
$ws.Range("P3").Value = 0.7077439387378263
$ws.Range("Q3").Value = 0.291790386
$ws.Range("R3").Value = 2.626113474
$ws.Range("S3").Value = 0.1749904784896267
$ws.Range("T3").Value = 0.1749904784896266

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Fgf13"
$ws.Range("C4").Value = "Scn8a"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.302256
$ws.Range("H4").Value = 0.906768
$ws.Range("I4").Value = 0.2472511157095891
$ws.Range("J4").Value = 0.2472511157095891
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.1417423333333333
$ws.Range("N4").Value = 0.425227
$ws.Range("O4").Value = 0.1039153461392964
$ws.Range("P4").Value = 0.1039153461392964
$ws.Range("Q4").Value = 0.042842470704
$ws.Range("R4").Value = 0.385582236336
$ws.Range("S4").Value = 0.02569318527228917
$ws.Range("T4").Value = 0.02569318527228917

$ws.Range("A5").Value = "Inflammatory-Mac"
$ws.Range("B5").Value = "Fgf13"
$ws.Range("C5").Value = "Scn8a"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.008244333333333334
$ws.Range("H5").Value = 0.024733
$ws.Range("I5").Value = 0.006744020350128443
$ws.Range("J5").Value = 0.006744020350128443
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.2569
$ws.Range("N5").Value = 0.7706999999999999
$ws.Range("O5").Value = 0.1883407151228772
$ws.Range("P5").Value = 0.1883407151228772
$ws.Range("Q5").Value = 0.002117969233333333
$ws.Range("R5").Value = 0.0190617231
$ws.Range("S5").Value = 0.001270173615546428
$ws.Range("T5").Value = 0.001270173615546427

$ws.Range("A6").Value = "Inflammatory-Mac"
$ws.Range("B6").Value = "Fgf13"
$ws.Range("C6").Value = "Scn8a"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.008244333333333334
$ws.Range("H6").Value = 0.024733
$ws.Range("I6").Value = 0.006744020350128443
$ws.Range("J6").Value = 0.006744020350128443
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.965375
$ws.Range("N6").Value = 2.896125
$ws.Range("O6").Value = 0.7077439387378264
$ws.Range("P6").Value = 0.7077439387378263
$ws.Range("Q6").Value = 0.007958873291666668
$ws.Range("R6").Value = 0.071629859625
$ws.Range("S6").Value = 0.004773039525527959
$ws.Range("T6").Value = 0.004773039525527958

$ws.Range("A7").Value = "Inflammatory-Mac"
$ws.Range("B7").Value = "Fgf13"
$ws.Range("C7").Value = "Scn8a"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.008244333333333334
$ws.Range("H7").Value = 0.024733
$ws.Range("I7").Value = 0.006744020350128443
$ws.Range("J7").Value = 0.006744020350128443
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.1417423333333333
$ws.Range("N7").Value = 0.425227
$ws.Range("O7").Value = 0.1039153461392964
$ws.Range("P7").Value = 0.1039153461392964
$ws.Range("Q7").Value = 0.001168571043444445
$ws.Range("R7").Value = 0.010517139391
$ws.Range("S7").Value = 0.0007008072090540559
$ws.Range("T7").Value = 0.0007008072090540558

$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "Fgf13"
$ws.Range("C8").Value = "Scn8a"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.592361
$ws.Range("H8").Value = 1.777083
$ws.Range("I8").Value = 0.4845624839634215
$ws.Range("J8").Value = 0.4845624839634214
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.2569
$ws.Range("N8").Value = 0.7706999999999999
$ws.Range("O8").Value = 0.1883407151228772
$ws.Range("P8").Value = 0.1883407151228772
$ws.Range("Q8").Value = 0.1521775409
$ws.Range("R8").Value = 1.3695978681
$ws.Range("S8").Value = 0.09126284475138852
$ws.Range("T8").Value = 0.0912628447513885

$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "Fgf13"
$ws.Range("C9").Value = "Scn8a"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.592361
$ws.Range("H9").Value = 1.777083
$ws.Range("I9").Value = 0.4845624839634215
$ws.Range("J9").Value = 0.4845624839634214
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.965375
$ws.Range("N9").Value = 2.896125
$ws.Range("O9").Value = 0.7077439387378264
$ws.Range("P9").Value = 0.7077439387378263
$ws.Range("Q9").Value = 0.5718505003750001
$ws.Range("R9").Value = 5.146654503375
$ws.Range("S9").Value = 0.3429461609648567
$ws.Range("T9").Value = 0.3429461609648566

$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Fgf13"
$ws.Range("C10").Value = "Scn8a"
$ws.Range("D10").Value = "MuSCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.592361
$ws.Range("H10").Value = 1.777083
$ws.Range("I10").Value = 0.4845624839634215
$ws.Range("J10").Value = 0.4845624839634214
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.1417423333333333
$ws.Range("N10").Value = 0.425227
$ws.Range("O10").Value = 0.1039153461392964
$ws.Range("P10").Value = 0.1039153461392964
$ws.Range("Q10").Value = 0.08396263031566667
$ws.Range("R10").Value = 0.755663672841
$ws.Range("S10").Value = 0.05035347824717619
$ws.Range("T10").Value = 0.05035347824717618

$ws.Range("A11").Value = "Resolving-Mac"
$ws.Range("B11").Value = "Fgf13"
$ws.Range("C11").Value = "Scn8a"
$ws.Range("D11").Value = "ECs"
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 0.6666666666666666
$ws.Range("G11").Value = 0.3196043333333333
$ws.Range("H11").Value = 0.9588129999999999
$ws.Range("I11").Value = 0.261442379976861
$ws.Range("J11").Value = 0.2614423799768609
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 0.6666666666666666
$ws.Range("M11").Value = 0.2569
$ws.Range("N11").Value = 0.7706999999999999
$ws.Range("O11").Value = 0.1883407151228772
$ws.Range("P11").Value = 0.1883407151228772
$ws.Range("Q11").Value = 0.08210635323333332
$ws.Range("R11").Value = 0.7389571790999999
$ws.Range("S11").Value = 0.04924024480826899
$ws.Range("T11").Value = 0.04924024480826897

$ws.Range("A12").Value = "Resolving-Mac"
$ws.Range("B12").Value = "Fgf13"
$ws.Range("C12").Value = "Scn8a"
$ws.Range("D12").Value = "FAPs"
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 0.6666666666666666
$ws.Range("G12").Value = 0.3196043333333333
$ws.Range("H12").Value = 0.9588129999999999
$ws.Range("I12").Value = 0.261442379976861
$ws.Range("J12").Value = 0.2614423799768609
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 0.965375
$ws.Range("N12").Value = 2.896125
$ws.Range("O12").Value = 0.7077439387378264
$ws.Range("P12").Value = 0.7077439387378263
$ws.Range("Q12").Value = 0.3085380332916666
$ws.Range("R12").Value = 2.776842299625
$ws.Range("S12").Value = 0.185034259757815
$ws.Range("T12").Value = 0.185034259757815

$ws.Range("A13").Value = "Resolving-Mac"
$ws.Range("B13").Value = "Fgf13"
$ws.Range("C13").Value = "Scn8a"
$ws.Range("D13").Value = "MuSCs"
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 0.6666666666666666
$ws.Range("G13").Value = 0.3196043333333333
$ws.Range("H13").Value = 0.9588129999999999
$ws.Range("I13").Value = 0.261442379976861
$ws.Range("J13").Value = 0.2614423799768609
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.1417423333333333
$ws.Range("N13").Value = 0.425227
$ws.Range("O13").Value = 0.1039153461392964
$ws.Range("P13").Value = 0.1039153461392964
$ws.Range("Q13").Value = 0.04530146395011111
$ws.Range("R13").Value = 0.407713175551
$ws.Range("S13").Value = 0.02716787541077696
$ws.Range("T13").Value = 0.02716787541077695

